$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing content fix: D3 question text gained inline math ---
$ws.Cells.Item(3, 4).Value = "Another question with \( R = a^i/se \) math?"

# --- New header cells K1:M1 ---
$ws.Cells.Item(1, 11).Value = "Nickname"
$ws.Cells.Item(1, 12).Value = "Human TEKS"
$ws.Cells.Item(1, 13).Value = "Machine TEKS"

# --- New data row 2 (K2:M2) ---
$ws.Cells.Item(2, 11).Value = "Q1"
$ws.Cells.Item(2, 12).Value = "T1.1"
$ws.Cells.Item(2, 13).Value = "dacf53a6-2b09-49f1-9926-de4efe1049e0"

# --- New data row 3 (K3:M3) ---
$ws.Cells.Item(3, 11).Value = "Q2"
$ws.Cells.Item(3, 12).Value = "T1.2"
$ws.Cells.Item(3, 13).Value = "c6623b8d-1eb7-41bf-875b-3456036000f9"

# --- Column widths ---
# Column D widened to fit the longer question text.
$ws.Columns.Item(4).ColumnWidth = 38
# New columns L (Human TEKS) and M (Machine TEKS) sized to fit their content.
$ws.Columns.Item(12).ColumnWidth = 10.666666666666666
$ws.Columns.Item(13).ColumnWidth = 34.833333333333336

# --- Selection / view ---
$ws.Range("D2").Select() | Out-Null

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1 | Out-Null
